# Update the lattice-multiplication exercise table: every cell's problem
# (title, top split-digits, dashes, and the two left split-digit rows) is
# replaced with a new set of values while keeping the same layout/format.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11

# row, col, title, topSplit, row1, row2
$cells = @(
    @(1, 1, "90 x 96", "  9    6", "9|    |", "0|    |"),
    @(1, 2, "54 x 58", "  5    8", "5|    |", "4|    |"),
    @(1, 3, "62 x 48", "  4    8", "6|    |", "2|    |"),
    @(2, 1, "79 x 13", "  1    3", "7|    |", "9|    |"),
    @(2, 2, "96 x 39", "  3    9", "9|    |", "6|    |"),
    @(2, 3, "92 x 14", "  1    4", "9|    |", "2|    |"),
    @(3, 1, "15 x 22", "  2    2", "1|    |", "5|    |"),
    @(3, 2, "45 x 64", "  6    4", "4|    |", "5|    |"),
    @(3, 3, "98 x 90", "  9    0", "9|    |", "8|    |"),
    @(4, 1, "14 x 74", "  7    4", "1|    |", "4|    |"),
    @(4, 2, "34 x 35", "  3    5", "3|    |", "4|    |"),
    @(4, 3, "73 x 70", "  7    0", "7|    |", "3|    |"),
    @(5, 1, "78 x 70", "  7    0", "7|    |", "8|    |"),
    @(5, 2, "46 x 30", "  3    0", "4|    |", "6|    |"),
    @(5, 3, "85 x 95", "  9    5", "8|    |", "5|    |")
)

foreach ($row in $cells) {
    $r = $row[0]
    $c = $row[1]
    $title = $row[2]
    $topSplit = $row[3]
    $line1 = $row[4]
    $line2 = $row[5]

    $newText = $title + $vt + $topSplit + $vt + "  ----" + $vt + $line1 + $vt + $line2
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $newText
}
